# The author's commit ("Added regular expressions slides") trimmed a
# leftover decorative shape from the "Nested Types and Anonymous
# Classes" deck: slide 21 had two near-identical small dark rectangles
# ("object 6" and "object 7") stacked side by side near the top of the
# slide; "object 7" is removed, leaving "object 6" untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)

$shape = $s.Shapes.Item("object 7")
$shape.Delete()
